$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (price) as text so numeric-looking strings like "0.999"
# or "49.06" are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "71.010.65"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.797.03"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "699.26"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "169.89"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").Value = "3.797.39"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "7.54"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").Value = "0.479"
$ws.Range("E12").Value = "  +4.08%  "
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "36.22"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").Value = "4.438.03"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "3.831.36"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "71.093.61"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "17.60"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "511.56"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").Value = "10.46"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "0.716"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").Value = "83.48"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("E25").Value = "  -3.54%  "
$ws.Range("D26").Value = "12.67"
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("D27").Value = "3.944.29"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -5.06%  "
$ws.Range("D31").Value = "2.96"
$ws.Range("E31").Value = "  -5.14%  "
$ws.Range("D32").Value = "2.27"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").Value = "7.30"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").Value = "29.11"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").Value = "0.172"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("D36").Value = "9.33"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "3.762.23"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "6.69"
$ws.Range("E39").Value = "  +11.11%  "
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D44").Value = "3.18"
$ws.Range("E44").Value = "  -5.85%  "
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "163.79"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D49").Value = "424.74"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "8.65"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "1.37"
$ws.Range("E51").Value = "  -1.04%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "0.000304"
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "49.06"
$ws.Range("E48").Value = "  +0.15%  "

# Restore the default (unstyled) cell style now that the text values are set,
# so the cells don't end up carrying an explicit style index.
$priceRange.Style = "Normal"
